$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 ("I0") and J1 ("IF") with the same style as existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-67: columns I (I0) and J (IF)
$data = @(
    @(2, 8, 8),
    @(3, 5, 5),
    @(4, 6, 6),
    @(5, 6, 6),
    @(6, 10, 10),
    @(7, 6, 6),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 10, 10),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 8, 8),
    @(17, 9, 10),
    @(18, 8, 8),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 11, 11),
    @(22, 8, 8),
    @(23, 9, 9),
    @(24, 9, 9),
    @(25, 9, 9),
    @(26, 9, 9),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 9, 9),
    @(31, 9, 9),
    @(32, 9, 9),
    @(33, 8, 8),
    @(34, 9, 9),
    @(35, 9, 9),
    @(36, 9, 9),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 10, 10),
    @(41, 9, 9),
    @(42, 9, 9),
    @(43, 9, 9),
    @(44, 9, 9),
    @(45, 9, 9),
    @(46, 8, 8),
    @(47, 9, 10),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 9, 9),
    @(51, 9, 9),
    @(52, 9, 9),
    @(53, 9, 9),
    @(54, 9, 9),
    @(55, 9, 9),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 9, 9),
    @(59, 9, 9),
    @(60, 9, 9),
    @(61, 9, 9),
    @(62, 9, 9),
    @(63, 9, 9),
    @(64, 9, 9),
    @(65, 9, 9),
    @(66, 7, 7),
    @(67, 4, 4)
)

foreach ($item in $data) {
    $r = $item[0]
    $i0 = $item[1]
    $iF = $item[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
